# Update per-year worksheets (sheets 1-8) with corrected lexical diversity values
# (count as proportion of all tokens, not just word tokens)
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item(1)
$ws.Cells.Item(2, 2).Value = 0.3444866920152092
$ws = $wb.Worksheets.Item(3)
$ws.Cells.Item(2, 2).Value = 0.9565217391304348
$ws.Cells.Item(3, 2).Value = 0.9375
$ws.Cells.Item(4, 2).Value = 0.9032258064516129
$ws.Cells.Item(5, 2).Value = 0.8205128205128205
$ws.Cells.Item(6, 2).Value = 0.8235294117647058
$ws.Cells.Item(8, 2).Value = 0.5814977973568282
$ws.Cells.Item(9, 2).Value = 0.9090909090909091
$ws.Cells.Item(10, 2).Value = 0.88
$ws.Cells.Item(11, 2).Value = 0.92
$ws = $wb.Worksheets.Item(4)
$ws.Cells.Item(2, 2).Value = 0.9090909090909091
$ws.Cells.Item(3, 2).Value = 0.9285714285714286
$ws.Cells.Item(4, 2).Value = 0.7216494845360825
$ws.Cells.Item(5, 2).Value = 0.9230769230769231
$ws.Cells.Item(6, 2).Value = 0.96
$ws.Cells.Item(7, 2).Value = 0.8428571428571429
$ws = $wb.Worksheets.Item(5)
$ws.Cells.Item(2, 2).Value = 0.9629629629629629
$ws.Cells.Item(4, 2).Value = 0.9047619047619048
$ws.Cells.Item(5, 2).Value = 0.8947368421052632
$ws.Cells.Item(6, 2).Value = 0.8947368421052632
$ws.Cells.Item(7, 2).Value = 0.84
$ws.Cells.Item(8, 2).Value = 0.8666666666666667
$ws.Cells.Item(9, 2).Value = 0.3605633802816902
$ws.Cells.Item(10, 2).Value = 0.3406593406593407
$ws.Cells.Item(12, 2).Value = 0.8214285714285714
$ws.Cells.Item(13, 2).Value = 0.8076923076923077
$ws = $wb.Worksheets.Item(6)
$ws.Cells.Item(2, 2).Value = 0.7763157894736842
$ws.Cells.Item(3, 2).Value = 0.6590909090909091
$ws.Cells.Item(4, 2).Value = 0.5955056179775281
$ws.Cells.Item(5, 2).Value = 0.9259259259259259
$ws.Cells.Item(6, 2).Value = 0.9411764705882353
$ws.Cells.Item(7, 2).Value = 0.9
$ws.Cells.Item(9, 2).Value = 0.9117647058823529
$ws.Cells.Item(10, 2).Value = 0.7714285714285715
$ws.Cells.Item(11, 2).Value = 0.8125
$ws.Cells.Item(12, 2).Value = 0.8888888888888888
$ws.Cells.Item(13, 2).Value = 0.8125
$ws.Cells.Item(14, 2).Value = 0.8235294117647058
$ws.Cells.Item(15, 2).Value = 0.7037037037037037
$ws.Cells.Item(16, 2).Value = 0.96
$ws.Cells.Item(17, 2).Value = 0.8709677419354839
$ws.Cells.Item(18, 2).Value = 0.9259259259259259
$ws.Cells.Item(19, 2).Value = 0.8387096774193549
$ws.Cells.Item(20, 2).Value = 0.9130434782608695
$ws.Cells.Item(21, 2).Value = 0.88
$ws.Cells.Item(22, 2).Value = 0.8157894736842105
$ws.Cells.Item(23, 2).Value = 0.9545454545454546
$ws.Cells.Item(24, 2).Value = 0.625
$ws.Cells.Item(26, 2).Value = 0.9285714285714286
$ws.Cells.Item(27, 2).Value = 0.8928571428571429
$ws.Cells.Item(28, 2).Value = 0.8076923076923077
$ws.Cells.Item(29, 2).Value = 0.6666666666666666
$ws.Cells.Item(30, 2).Value = 0.9130434782608695
$ws = $wb.Worksheets.Item(7)
$ws.Cells.Item(2, 2).Value = 0.8571428571428571
$ws.Cells.Item(3, 2).Value = 0.85
$ws.Cells.Item(5, 2).Value = 0.9545454545454546
$ws.Cells.Item(7, 2).Value = 0.9047619047619048
$ws.Cells.Item(8, 2).Value = 0.9615384615384616
$ws.Cells.Item(9, 2).Value = 0.72
$ws.Cells.Item(10, 2).Value = 0.7746478873239436
$ws.Cells.Item(11, 2).Value = 0.9629629629629629
$ws.Cells.Item(12, 2).Value = 0.9565217391304348
$ws.Cells.Item(13, 2).Value = 0.8
$ws.Cells.Item(15, 2).Value = 0.7777777777777778
$ws.Cells.Item(16, 2).Value = 0.9473684210526315
$ws.Cells.Item(18, 2).Value = 0.8461538461538461
$ws.Cells.Item(19, 2).Value = 0.8181818181818182
$ws.Cells.Item(20, 2).Value = 0.7142857142857143
$ws.Cells.Item(21, 2).Value = 0.8695652173913043
$ws.Cells.Item(22, 2).Value = 0.384995064165844
$ws.Cells.Item(23, 2).Value = 0.8461538461538461
$ws.Cells.Item(24, 2).Value = 0.9583333333333334
$ws.Cells.Item(25, 2).Value = 0.9629629629629629
$ws.Cells.Item(26, 2).Value = 0.8928571428571429
$ws.Cells.Item(27, 2).Value = 0.9230769230769231
$ws.Cells.Item(28, 2).Value = 0.9473684210526315
$ws.Cells.Item(29, 2).Value = 0.9615384615384616
$ws.Cells.Item(30, 2).Value = 0.9130434782608695
$ws.Cells.Item(31, 2).Value = 0.8545454545454545
$ws.Cells.Item(32, 2).Value = 0.9166666666666666
$ws.Cells.Item(33, 2).Value = 0.8064516129032258
$ws.Cells.Item(35, 2).Value = 0.8461538461538461
$ws.Cells.Item(36, 2).Value = 0.5158730158730159
$ws.Cells.Item(38, 2).Value = 0.8695652173913043
$ws.Cells.Item(39, 2).Value = 0.8666666666666667
$ws.Cells.Item(40, 2).Value = 0.8333333333333334
$ws.Cells.Item(41, 2).Value = 0.96
$ws.Cells.Item(42, 2).Value = 0.88
$ws.Cells.Item(43, 2).Value = 0.9642857142857143
$ws.Cells.Item(44, 2).Value = 0.8823529411764706
$ws.Cells.Item(45, 2).Value = 0.9166666666666666
$ws.Cells.Item(46, 2).Value = 0.5340136054421769
$ws.Cells.Item(47, 2).Value = 0.8064516129032258
$ws.Cells.Item(48, 2).Value = 0.9130434782608695
$ws.Cells.Item(50, 2).Value = 0.8461538461538461
$ws.Cells.Item(51, 2).Value = 0.8275862068965517
$ws.Cells.Item(52, 2).Value = 0.9333333333333333
$ws.Cells.Item(53, 2).Value = 0.8928571428571429
$ws.Cells.Item(54, 2).Value = 0.7068965517241379
$ws.Cells.Item(55, 2).Value = 0.8333333333333334
$ws.Cells.Item(56, 2).Value = 0.96
$ws.Cells.Item(57, 2).Value = 0.8
$ws.Cells.Item(58, 2).Value = 0.7428571428571429
$ws.Cells.Item(59, 2).Value = 0.90625
$ws.Cells.Item(60, 2).Value = 0.896551724137931
$ws.Cells.Item(61, 2).Value = 0.7857142857142857
$ws.Cells.Item(62, 2).Value = 0.9629629629629629
$ws.Cells.Item(63, 2).Value = 0.8275862068965517
$ws.Cells.Item(64, 2).Value = 0.9629629629629629
$ws.Cells.Item(65, 2).Value = 0.9523809523809523
$ws.Cells.Item(66, 2).Value = 0.8571428571428571
$ws = $wb.Worksheets.Item(8)
$ws.Cells.Item(2, 2).Value = 0.8837209302325582
$ws.Cells.Item(3, 2).Value = 0.75
$ws.Cells.Item(4, 2).Value = 0.8787878787878788
$ws.Cells.Item(5, 2).Value = 0.7547169811320755
$ws.Cells.Item(6, 2).Value = 0.8
$ws.Cells.Item(7, 2).Value = 0.8085106382978723
$ws.Cells.Item(8, 2).Value = 0.78
$ws.Cells.Item(9, 2).Value = 0.8653846153846154
$ws.Cells.Item(10, 2).Value = 0.8113207547169812
$ws.Cells.Item(11, 2).Value = 0.803921568627451
$ws.Cells.Item(12, 2).Value = 0.8823529411764706
$ws.Cells.Item(13, 2).Value = 0.8301886792452831
$ws.Cells.Item(14, 2).Value = 0.8301886792452831
$ws.Cells.Item(15, 2).Value = 0.7666666666666667
$ws.Cells.Item(16, 2).Value = 0.7454545454545455
$ws.Cells.Item(17, 2).Value = 0.8695652173913043
$ws.Cells.Item(18, 2).Value = 0.6349206349206349
$ws.Cells.Item(20, 2).Value = 0.7894736842105263
$ws.Cells.Item(21, 2).Value = 0.8032786885245902
$ws.Cells.Item(22, 2).Value = 0.8448275862068966
$ws.Cells.Item(23, 2).Value = 0.6862745098039216
$ws.Cells.Item(24, 2).Value = 0.8157894736842105
$ws.Cells.Item(25, 2).Value = 0.65625
$ws.Cells.Item(26, 2).Value = 0.7105263157894737
$ws.Cells.Item(27, 2).Value = 0.7741935483870968
$ws.Cells.Item(28, 2).Value = 0.86
$ws.Cells.Item(29, 2).Value = 0.7959183673469388
$ws.Cells.Item(30, 2).Value = 0.7551020408163265
$ws.Cells.Item(31, 2).Value = 0.7333333333333333

# Update Summary worksheet (sheet 9) with recomputed describe() stats
$ws = $wb.Worksheets.Item(9)
$ws.Cells.Item(2, 3).Value = 0.3444866920152092
$ws.Cells.Item(2, 5).Value = 0.3444866920152092
$ws.Cells.Item(2, 6).Value = 0.3444866920152092
$ws.Cells.Item(2, 7).Value = 0.3444866920152092
$ws.Cells.Item(2, 8).Value = 0.3444866920152092
$ws.Cells.Item(2, 9).Value = 0.3444866920152092
$ws.Cells.Item(3, 3).Value = 1
$ws.Cells.Item(3, 5).Value = 1
$ws.Cells.Item(3, 6).Value = 1
$ws.Cells.Item(3, 7).Value = 1
$ws.Cells.Item(3, 8).Value = 1
$ws.Cells.Item(3, 9).Value = 1
$ws.Cells.Item(4, 3).Value = 0.8731878484307311
$ws.Cells.Item(4, 4).Value = 0.1163876226659042
$ws.Cells.Item(4, 5).Value = 0.5814977973568282
$ws.Cells.Item(4, 6).Value = 0.8376470588235294
$ws.Cells.Item(4, 7).Value = 0.906158357771261
$ws.Cells.Item(4, 8).Value = 0.933125
$ws.Cells.Item(4, 9).Value = 1
$ws.Cells.Item(5, 3).Value = 0.8808743146887476
$ws.Cells.Item(5, 4).Value = 0.0870758026025065
$ws.Cells.Item(5, 5).Value = 0.7216494845360825
$ws.Cells.Item(5, 6).Value = 0.8594155844155844
$ws.Cells.Item(5, 7).Value = 0.916083916083916
$ws.Cells.Item(5, 8).Value = 0.9271978021978022
$ws.Cells.Item(5, 9).Value = 0.96
$ws.Cells.Item(6, 3).Value = 0.8078507348886642
$ws.Cells.Item(6, 4).Value = 0.2226705062357575
$ws.Cells.Item(6, 5).Value = 0.3406593406593407
$ws.Cells.Item(6, 6).Value = 0.8179945054945055
$ws.Cells.Item(6, 7).Value = 0.880701754385965
$ws.Cells.Item(6, 8).Value = 0.9193121693121693
$ws.Cells.Item(6, 9).Value = 1
$ws.Cells.Item(7, 3).Value = 0.8505047590181406
$ws.Cells.Item(7, 4).Value = 0.1115018123473218
$ws.Cells.Item(7, 5).Value = 0.5955056179775281
$ws.Cells.Item(7, 6).Value = 0.8088942307692308
$ws.Cells.Item(7, 7).Value = 0.8844444444444444
$ws.Cells.Item(7, 8).Value = 0.9259259259259259
$ws.Cells.Item(7, 9).Value = 1
$ws.Cells.Item(8, 3).Value = 0.8702526621119039
$ws.Cells.Item(8, 4).Value = 0.116804107716692
$ws.Cells.Item(8, 5).Value = 0.384995064165844
$ws.Cells.Item(8, 6).Value = 0.8275862068965517
$ws.Cells.Item(8, 7).Value = 0.8928571428571429
$ws.Cells.Item(8, 8).Value = 0.9583333333333334
$ws.Cells.Item(8, 9).Value = 1
$ws.Cells.Item(9, 3).Value = 0.7973556093127645
$ws.Cells.Item(9, 4).Value = 0.07473090325622489
$ws.Cells.Item(9, 5).Value = 0.6349206349206349
$ws.Cells.Item(9, 6).Value = 0.7548132460531383
$ws.Cells.Item(9, 7).Value = 0.8016393442622951
$ws.Cells.Item(9, 8).Value = 0.8411678594664932
$ws.Cells.Item(9, 9).Value = 1
